# External Antenna selected and buddy LED blinking when comm. started
#
# - Update the "Date:" header (C1) to the new date serial.
# - Antenna (row 9) status: Open -> Closed.
# - "Improve streaming performance" (row 45) status: Open -> Ongoing.
# - "VU meter" (row 52) status: Open -> Rejected, with a note in D52
#   explaining why ("Using LED blink instead").
# - "Add WiFi signal strength to cloud API and app" (row 54) status:
#   Ongoing -> Closed.
# - Move the viewport/selection to the area that was just edited.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Date: header cell (serial date number, keeps existing date format)
$ws.Range("C1").Value = 43119

# Antenna -> Closed
$ws.Range("C9").Value = "Closed"

# Improve streaming performance -> Ongoing
$ws.Range("C45").Value = "Ongoing"

# VU meter -> Rejected, with explanation note
$ws.Range("C52").Value = "Rejected"
$ws.Range("D52").Value = "Using LED blink instead"

# Add WiFi signal strength to cloud API and app -> Closed
$ws.Range("C54").Value = "Closed"

# Reflect the place the author was last looking at / editing
$ws.Range("D52").Select()
